# Fruta / hortaliza, semanal
# Insert a new weekly record above row 377 (pushing existing rows 377:431 down to 378:432)
# and populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 377; this shifts rows 377:431 down to 378:432
# and automatically extends the sheet dimension to A1:T432.
$ws.Rows("377:377").Insert()

# Populate the newly inserted row 377 with the new observation.
$ws.Range("A377").Value = 10
$ws.Range("B377").Value = "Vega Modelo de Temuco"
$ws.Range("C377").Value = "La Araucanía"
$ws.Range("D377").Value = 45077
$ws.Range("E377").Value = 9
$ws.Range("F377").Value = "Fruta"
$ws.Range("G377").Value = 100102
$ws.Range("H377").Value = "Cítricos"
$ws.Range("I377").Value = 100102006
$ws.Range("J377").Value = "Pomelo"
$ws.Range("K377").Value = "Start Ruby"
$ws.Range("L377").Value = "Primera"
$ws.Range("M377").Value = 80
$ws.Range("N377").Value = 15000
$ws.Range("O377").Value = 15000
$ws.Range("P377").Value = 15000
$ws.Range("Q377").Value = "$/bandeja 15 kilos granel"
$ws.Range("R377").Value = "Región de O'Higgins"
$ws.Range("S377").Value = 1000
$ws.Range("T377").Value = 15
